# Re-order / update rows on the "BlewerResume.docx" and
# "RinglerShawn_Resume.docx" sheets to match the refreshed scrape data.
# (rows were reshuffled and a handful of metric columns + contact-email
# lists were refreshed; row counts per sheet are unchanged.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BlewerResume.docx")
$ws.Range("A3").Value = "Google plus"
$ws.Range("B3").Value = "https://raygun.com/blog/programming-languages/"
$ws.Range("C3").Value = 127
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = "example@xyz.com`n"
$ws.Range("A4").Value = "Access to this page has been denied."
$ws.Range("B4").Value = "https://www.upwork.com/l/cn/python-developers/"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").ClearContents() | Out-Null
$ws.Range("A10").Value = "Remote C++ Jobs in December 2019"
$ws.Range("B10").Value = "https://remoteok.io/remote-c-plus-plus-jobs"
$ws.Range("C10").Value = 11167
$ws.Range("D10").Value = 79
$ws.Range("E10").Value = 767
$ws.Range("F10").Value = "hello@geektastic.com`njobs@alienskin.com`nschuss@madeinoffice.com`njobs@komodoplatform.com`n"
$ws.Range("A11").Value = "Remote Software Developer Jobs in December 2019"
$ws.Range("B11").Value = "https://remoteok.io/remote-dev-jobs"
$ws.Range("C11").Value = 10880
$ws.Range("D11").Value = 94
$ws.Range("E11").Value = 1534
$ws.Range("F11").Value = "hello@geektastic.com`ntalent@tuftandneedle.com`nsentinelhr@sentinel.com`njobs@wpwhitesecurity.com`nalise.moncure@integratedrental.com`nalise.moncure@integratedrental.com`n"
$ws.Range("A12").Value = "GitHub - uhub/awesome-cpp: A curated list of awesome C++ frameworks, libraries and software."
$ws.Range("B12").Value = "https://github.com/uhub/awesome-cpp"
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 19
$ws.Range("E12").Value = 49
$ws.Range("F12").Value = "sean@seanstarkey.com`nsean@seanstarkey.com`n"
$ws.Range("A13").Value = "Ask HN: Who wants to be hired? (April 2019) | Hacker News"
$ws.Range("B13").Value = "https://news.ycombinator.com/item?id=19543938"
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = 26
$ws.Range("E13").Value = 323
$ws.Range("F13").Value = "nathompson7@protonmail.com`npcombs@gmail.com`nishtaygrud@hotmail.com`ncraig.glennie@gmail.com`nnbp@fynestro.com`nhi@petekeen.net`nmail@alekseilevin.com`nnithya1810@gmail.com`nlucaschaufournier@gmail.com`nclaudio.salinitro@gmail.com`nforjob18846@gmail.com`nhnjobs@gmail.com`nrobby.ronk@gmail.com`nsayhar@gmail.com`nadityanalluri9@gmail.com`nserge@redmondanalytics.com`nluisenriquenovoa@gmail.com`ntyteen4a03@gmail.com`nbiz@harlanji.com`nhn@gmail.com`nalbertommoura@gmail.com`nbaile320@umn.edu`nriveracarvelli@gmail.com`nlogan@logankoester.com`nhn@lazerwalker.com`nroee@sandsquid.com`ndustin.freeman@gmail.com`nunleashit@protonmail.com`ntaras.brizitsky@gmail.com`nnityamd@gmail.com`nabdurleo91@gmail.com`nhervan@gmail.com`na31415926535@gmail.com`nhginfla@gmail.com`npaul.wujek@gmail.com`nsyed@shuttari.com`njoseph@jwdougherty.com`neg.public@gmail.com`nshashanksira@gmail.com`nrivera.julioa@gmail.com`nbartus.csongor@gmail.com`nmarcospereiradev@gmail.com`nchristian.kopac@protonmail.com`nsusmus7@gmail.com`ntiwatson@gmail.com`nabdelgzali@gmail.com`nmichail@rybakov.com`nwlodzislav@outlook.com`nbyshinyo@gmail.com`nhireme@gmail.com`njon@dontbreakthebuild.com`nfranklai@protonmail.com`nrazzintown@gmail.com`njesse.hughes.it@gmail.com`njrudisill@gmail.com`nhusseinlotfizaki@gmail.com`nsofikurochkina@gmail.com`ncch5ng.job@gmail.com`nhackernews@ronilan.com`nandre.b.mourao@gmail.com`nbjpcjp@gmail.com`nchancelor.oneal@gmail.com`nchevonied@gmail.com`nbradleydavidhoffman@gmail.com`nhi@mrassili.com`njustin.crandell.developer@gmail.com`njoe@bokengroup.com`nakshay.sharma09695@gmail.com`ncymenvig@gmail.com`nmichael.nicolaou@protonmail.com`nshivrajnesargi07@gmail.com`nshubh065@gmail.com`nsimon.borer@gmail.com`ndixiekorley@gmail.com`nchris.comeau@skyriser.com`niodbh@iodbh.net`nfullstackguy@gmail.com`nhtavenido@gmail.com`nadige01can@gmail.com`nrichardbryancall@gmail.com`nbenjamindhsu@gmail.com`nlars.jarlvik@gmail.com`ntucker.r.chapman@gmail.com`nhello@jensunltd.com`nat@gmail.com`njialunz@umich.edu`nqoutland@gmail.com`nthibjp@protonmail.com`ncyrilbenson47@gmail.com`nchc4000@gmail.com`nhnhiringthread@gmail.com`ndsm140130@utdallas.edu`narawde@gmail.com`nrrrasti@yahoo.com`nmmchristian@gmail.com`natiwari3@binghamton.edu`ntoresson.gustav@gmail.com`nemile.senga@gmail.com`n369@holbertonschool.com`nmax42@protonmail.com`narturo@volpe.com`ntekadeaniruddha@gmail.com`ncharles@techascent.com`nerudyn@protonmail.com`nrajesh.singh@vanderbilt.edu`nmlin36729@gmail.com`ndsuryav@gmail.com`nroscoe1245@gmail.com`nhn@geoffgolder.com`nedimaudo@gmail.com`nbrandons.519@gmail.com`nrovilayjnr01@gmail.com`nsrtplayer@gmail.com`nadonismurati@gmail.com`nisharamet@gmail.com`nbirykovegor@gmail.com`nvishnusharathr@gmail.com`nbxeighty8@gmail.com`nbxeighty8@gmail.com`nkhan.alistar@gmail.com`nclintfidel@gmail.com`nbrian@angularjobs.com`ncannadayr@gmail.com`nkellybhoward@gmail.com`nandrey.boar@gmail.com`nkeaneswolter@gmail.com`nheyarviind@gmail.com`nhn@ycombinator.com`nhn@ycombinator.com`n"
$ws.Range("A14").Value = "Should C# or C++ be chosen for learning Games Programming (consoles)? - Stack Overflow"
$ws.Range("B14").Value = "https://stackoverflow.com/questions/2203093/should-c-sharp-or-c-be-chosen-for-learning-games-programming-consoles"
$ws.Range("C14").Value = 62
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 14
$ws.Range("A15").Value = "What is the difference between C, C++ and C#? - Quora"
$ws.Range("B15").Value = "https://www.quora.com/What-is-the-difference-between-C-C++-and-C"
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = 5
$ws.Range("A16").Value = "C++ Rest Sdk Linux"
$ws.Range("B16").Value = "https://rtmm.ts-fliesenservice.de/c++-rest-sdk-linux.html"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 8
$ws.Range("F16").ClearContents() | Out-Null
$ws.Range("A17").Value = "Which is better for programming, a Mac or Windows laptop? - Quora"
$ws.Range("B17").Value = "https://www.quora.com/Which-is-better-for-programming-a-Mac-or-Windows-laptop"
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 16

$ws = $wb.Worksheets.Item("RinglerShawn_Resume.docx")
$ws.Range("A3").Value = "Warning: Your programming career - SoloLearn - Medium"
$ws.Range("B3").Value = "https://medium.com/sololearn/warning-your-programming-career-b9579b3a878b"
$ws.Range("C3").Value = 54
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 20
$ws.Range("A4").Value = "The 7 Most In-Demand Programming Languages of 2019 - Coding Dojo Blog"
$ws.Range("B4").Value = "https://www.codingdojo.com/blog/the-7-most-in-demand-programming-languages-of-2019"
$ws.Range("C4").Value = 59
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 11
$ws.Range("A5").Value = "Which Programming Language Should You Learn Next?"
$ws.Range("B5").Value = "https://www.freecodecamp.org/news/which-programming-language-should-you-learn-next-487d077baa32/"
$ws.Range("C5").Value = 13
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 9
$ws.Range("A6").Value = "Top 10 Programming Languages of the World – 2019 to begin with… - GeeksforGeeks"
$ws.Range("B6").Value = "https://www.geeksforgeeks.org/top-10-programming-languages-of-the-world-2019-to-begin-with/"
$ws.Range("C6").Value = 28
$ws.Range("D6").Value = 41
$ws.Range("E6").Value = 29
$ws.Range("F6").ClearContents() | Out-Null
$ws.Range("A7").Value = "Thinkful Logo Black@2x"
$ws.Range("B7").Value = "https://www.thinkful.com/blog/what-programming-language-should-you-learn-according-to-your-state/"
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "tatiana@thinkful.com`ntatiana@thinkful.com`nyou@example.com`nhello@thinkful.com`nsupport@thinkful.com`nsupport@thinkful.com`n"
$ws.Range("A8").Value = "What Programming Language Should a Beginner Learn in 2019? | Codementor"
$ws.Range("B8").Value = "https://www.codementor.io/codementorteam/beginner-programming-language-job-salary-community-7s26wmbm6"
$ws.Range("C8").Value = 68
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 24
$ws.Range("A13").Value = "Google plus"
$ws.Range("B13").Value = "https://raygun.com/blog/programming-languages/"
$ws.Range("C13").Value = 127
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = "example@xyz.com`n"
$ws.Range("A14").Value = "Javascript C++ Jobs, Employment | Indeed.com"
$ws.Range("B14").Value = "https://www.indeed.com/q-Javascript-C++-jobs.html"
$ws.Range("C14").Value = 1595
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 20
$ws.Range("F14").ClearContents() | Out-Null
$ws.Range("C15").Value = 1537
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 15
$ws.Range("A16").Value = "The 9 Best Programming Languages to Learn in 2019 | Fullstack Academy"
$ws.Range("B16").Value = "https://www.fullstackacademy.com/blog/nine-best-programming-languages-to-learn-2018"
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "hello@fullstackacademy.com`n"
$ws.Range("A17").Value = "The Best Programming Languages for each Situation"
$ws.Range("B17").Value = "https://tomassetti.me/best-programming-languages/"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 34
$ws.Range("F17").ClearContents() | Out-Null
